$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Kết quả import câu hỏi"
$ws.Range("B2").Value = "Mon Dec 18 2023 16:59:31 GMT+0700 (Indochina Time)"
$ws.Range("A3").Value = "Tổng câu hỏi"
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 1
